$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the progress notes for the "e2e, service order creating (map)" task
# in the "Date" rows 5 (41754) and 6 (41755), column E ("Толя" / task owner column).
$ws.Range("E5").Value = "e2e, service order creating (map)- 40%"
$ws.Range("E5").WrapText = $true

$ws.Range("E6").Value = "e2e, service order creating (map)- 60%"
$ws.Range("E6").WrapText = $true

# Reflect where the author ended up working: sheet scrolled so column C is
# the left-most visible column, with the newly-entered cell E6 selected.
$ws.Activate() | Out-Null
$ws.Range("E6").Select() | Out-Null
